{"js": "// Replace the 25 two-digit-division answer strings in the table with\n// their updated values (commit c8c62b6).\nconst replacements = [\n  [\"60\u00f74=15, 0\", \"88\u00f78=11, 0\"],\n  [\"29\u00f73=9, 2\", \"73\u00f75=14, 3\"],\n  [\"76\u00f78=9, 4\", \"27\u00f78=3, 3\"],\n  [\"84\u00f73=28, 0\", \"22\u00f72=11, 0\"],\n  [\"64\u00f72=32, 0\", \"18\u00f74=4, 2\"],\n  [\"24\u00f78=3, 0\", \"37\u00f75=7, 2\"],\n  [\"10\u00f78=1, 2\", \"59\u00f77=8, 3\"],\n  [\"76\u00f73=25, 1\", \"82\u00f76=13, 4\"],\n  [\"93\u00f75=18, 3\", \"55\u00f76=9, 1\"],\n  [\"37\u00f77=5, 2\", \"14\u00f73=4, 2\"],\n  [\"30\u00f79=3, 3\", \"11\u00f73=3, 2\"],\n  [\"46\u00f75=9, 1\", \"42\u00f75=8, 2\"],\n  [\"46\u00f73=15, 1\", \"52\u00f77=7, 3\"],\n  [\"50\u00f79=5, 5\", \"27\u00f74=6, 3\"],\n  [\"76\u00f79=8, 4\", \"79\u00f77=11, 2\"],\n  [\"61\u00f75=12, 1\", \"63\u00f72=31, 1\"],\n  [\"14\u00f72=7, 0\", \"21\u00f78=2, 5\"],\n  [\"53\u00f78=6, 5\", \"64\u00f74=16, 0\"],\n  [\"31\u00f73=10, 1\", \"94\u00f72=47, 0\"],\n  [\"39\u00f79=4, 3\", \"71\u00f77=10, 1\"],\n  [\"13\u00f72=6, 1\", \"68\u00f73=22, 2\"],\n  [\"28\u00f79=3, 1\", \"94\u00f79=10, 4\"],\n  [\"23\u00f74=5, 3\", \"91\u00f72=45, 1\"],\n  [\"23\u00f75=4, 3\", \"61\u00f74=15, 1\"],\n  [\"59\u00f78=7, 3\", \"19\u00f74=4, 3\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 two-digit-division answer strings in the table with\n# their updated values (commit c8c62b6).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"60\u00f74=15, 0\", \"88\u00f78=11, 0\"),\n    @(\"29\u00f73=9, 2\", \"73\u00f75=14, 3\"),\n    @(\"76\u00f78=9, 4\", \"27\u00f78=3, 3\"),\n    @(\"84\u00f73=28, 0\", \"22\u00f72=11, 0\"),\n    @(\"64\u00f72=32, 0\", \"18\u00f74=4, 2\"),\n    @(\"24\u00f78=3, 0\", \"37\u00f75=7, 2\"),\n    @(\"10\u00f78=1, 2\", \"59\u00f77=8, 3\"),\n    @(\"76\u00f73=25, 1\", \"82\u00f76=13, 4\"),\n    @(\"93\u00f75=18, 3\", \"55\u00f76=9, 1\"),\n    @(\"37\u00f77=5, 2\", \"14\u00f73=4, 2\"),\n    @(\"30\u00f79=3, 3\", \"11\u00f73=3, 2\"),\n    @(\"46\u00f75=9, 1\", \"42\u00f75=8, 2\"),\n    @(\"46\u00f73=15, 1\", \"52\u00f77=7, 3\"),\n    @(\"50\u00f79=5, 5\", \"27\u00f74=6, 3\"),\n    @(\"76\u00f79=8, 4\", \"79\u00f77=11, 2\"),\n    @(\"61\u00f75=12, 1\", \"63\u00f72=31, 1\"),\n    @(\"14\u00f72=7, 0\", \"21\u00f78=2, 5\"),\n    @(\"53\u00f78=6, 5\", \"64\u00f74=16, 0\"),\n    @(\"31\u00f73=10, 1\", \"94\u00f72=47, 0\"),\n    @(\"39\u00f79=4, 3\", \"71\u00f77=10, 1\"),\n    @(\"13\u00f72=6, 1\", \"68\u00f73=22, 2\"),\n    @(\"28\u00f79=3, 1\", \"94\u00f79=10, 4\"),\n    @(\"23\u00f74=5, 3\", \"91\u00f72=45, 1\"),\n    @(\"23\u00f75=4, 3\", \"61\u00f74=15, 1\"),\n    @(\"59\u00f78=7, 3\", \"19\u00f74=4, 3\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $found = $find.Execute(\n        [ref]$oldText, [ref]$false, [ref]$false, [ref]$false, [ref]$false,\n        [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$newText, [ref]2\n    )\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}\n"}
